# Update the "K" column (column G) values on Sheet1.
# These values were regenerated upstream (K instead of Strike#, recalculated
# std/mean, etc.) so we simply overwrite the literal numbers that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 3
    6  = 1
    7  = 3
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 1
    17 = 2
    18 = 0
    19 = 1
    20 = 3
    21 = 1
    22 = 0
    23 = 1
    24 = 3
    25 = 1
    26 = 1
    27 = 2
    28 = 0
    29 = 0
    30 = 0
    31 = 1
    32 = 1
    33 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
